$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 63 ("Hortaliza, Vega Modelo de
# Temuco - Arveja Verde" weekly refresh adds a new observation), pushing the
# former rows 63-81 down to 64-82.
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new weekly record.
$ws.Range("A63").Value = 10
$ws.Range("B63").Value = "Vega Modelo de Temuco"
$ws.Range("C63").Value = "La Araucanía"

# Date column keeps the same number format used by the rest of column D.
$ws.Range("D63").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D63").Value = 44559

$ws.Range("E63").Value = 9
$ws.Range("F63").Value = 100112022
$ws.Range("G63").Value = "Arveja Verde"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 110
$ws.Range("K63").Value = 15000
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = 15000
$ws.Range("N63").Value = "$/saco 25 kilos"
$ws.Range("O63").Value = "Región de La Araucanía"
$ws.Range("P63").Value = 600
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
